# Update the food table: translate category columns (B, D) from English to German,
# change the effort column (C) from a numeric scale to a German categorical label,
# and append a new "Tacos" row (row 16).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@('food', 'salty', 'effort', 'takeaway')
    ,@('Tortellini mit Käse Sahne Soße', 'herzhaft', 'bestellen', 'kochen')
    ,@('Nutella Brot', 'süß', 'wenig', 'kochen')
    ,@('Pizza bestellen', 'herzhaft', 'bestellen', 'bestellen')
    ,@('Vegetarisches Sushi', 'herzhaft', 'bestellen', 'bestellen')
    ,@('Ein Spiegelei', 'herzhaft', 'wenig', 'kochen')
    ,@('Tomate Mozarella Salat', 'herzhaft', 'wenig', 'kochen')
    ,@('Käsebrot', 'herzhaft', 'wenig', 'kochen')
    ,@('Sushi selber machen', 'herzhaft', 'hoch', 'kochen')
    ,@('Pizza selber machen', 'herzhaft', 'hoch', 'kochen')
    ,@('Nougat Bit Müsli', 'süß', 'wenig', 'kochen')
    ,@('Porridge mit Peanutbutter', 'süß', 'wenig', 'kochen')
    ,@('Cookies backen', 'süß', 'mittel', 'kochen')
    ,@('Hefezopf mit Nutella', 'süß', 'mittel', 'kochen')
    ,@('Obazda Brot', 'herzhaft', 'wenig', 'kochen')
    ,@('Tacos', 'herzhaft', 'mittel', 'kochen')
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 1
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

# Header row gets an explicit row height in the saved workbook.
$ws.Rows(1).RowHeight = 15

# The author left the cursor on F8 before saving.
$ws.Range("F8").Select()
